$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "P6_EL_DoubleEL" to "P6_EL_DualEL" in the two P6 rows (rows 24/25
# before the new row is inserted below; they become rows 25/26 afterward).
$ws.Range("E24").Value = "P6_EL_DualEL"
$ws.Range("E25").Value = "P6_EL_DualEL"

# Insert a new row above the current row 22 (P5 pre-toll-calibration run row),
# which copies formatting from the row above (row 21) matching the source
# workbook's pattern for this block of rows.
$ws.Rows("22:22").Insert()

# Populate the newly inserted row 22 with the Path4_02_pretollcalib run info.
$ws.Range("A22").Value = "NextGenFwys"
$ws.Range("B22").Value = 2035
$ws.Range("C22").Value = "2035_TM160_NGFr2_NP04_Path4_02_pretollcalib"
$ws.Range("D22").Value = "NGF_Round2"
$ws.Range("E22").Value = "P4_EL_PBA2050scope"
$ws.Range("F22").Value = "P4 pre-toll-calibration run"
$ws.Range("H22").Value = "NGF_Networks_NGFround2_P4_09"
$ws.Range("I22").Value = "https://app.asana.com/0/572982923864207/1207634640713913/f"
$ws.Range("J22").Value = "na"
$ws.Range("K22").Value = "na"

# Move the active selection to A22, matching the edited workbook's saved
# cursor position.
$ws.Range("A22").Select()
